{"js": "// Apply the same content edits Word made when the author revised the\n// chapter-9 example: (1) collapse the run-fragmented example title into a\n// single run of text, and (2) rename the R variables p_0 -> p0 and\n// se_0 -> se0 throughout the source-code listing.\n\n// 1) Merge the title runs \"Example \" + \"6\" + \":  \" + \"Car Stopping\" + \" \u2013\"\n//    + \" Significance Test About \" into one run with identical formatting.\n//    Re-inserting the exact same text (Replace) over the whole matched\n//    range makes the host collapse the formatting-identical runs into a\n//    single run, exactly like Word does when you retype/replace selected\n//    text.\nconst titleResults = context.document.body.search(\n  \"Example 6:  Car Stopping \\u2013 Significance Test About \",\n  { matchCase: true }\n);\ntitleResults.load(\"items\");\nawait context.sync();\n\nif (titleResults.items.length > 0) {\n  titleResults.items[0].insertText(\n    \"Example 6:  Car Stopping \\u2013 Significance Test About \",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// 2) Rename the R variables inside the annotated source-code block.\n//    Each occurrence lives entirely inside its own run, so a plain\n//    search + replace keeps every other run (and its formatting) intact.\nasync function replaceAll(findText, replaceText) {\n  const results = context.document.body.search(findText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\nawait replaceAll(\"p_0\", \"p0\");\nawait replaceAll(\"se_0\", \"se0\");\n", "ps1": "# Apply the same content edits Word made when the author revised the\n# chapter-9 example: (1) collapse the run-fragmented example title into a\n# single run of text, and (2) rename the R variables p_0 -> p0 and\n# se_0 -> se0 throughout the source-code listing.\n\n$d = $word.ActiveDocument\n\n# 1) Merge the title runs \"Example \" + \"6\" + \":  \" + \"Car Stopping\" + \" \u2013\"\n#    + \" Significance Test About \" into one run with identical formatting.\n#    Finding & replacing the exact same text over the whole matched range\n#    makes Word collapse the formatting-identical runs into a single run,\n#    exactly like retyping/replacing the selected text in the UI.\n$titleText = \"Example 6:  Car Stopping \" + [char]0x2013 + \" Significance Test About \"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $titleText\n$find.Replacement.Text = $titleText\n$find.MatchCase = $true\n# NOTE: deliberately leave the \"Wrap\" (11th) Execute argument as $null \u2014\n# passing a boxed boolean there trips this host's Find.Execute shim into\n# literally typing \"True\"/\"False\" into the document instead of toggling\n# wrap-around search.\n$find.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2) | Out-Null\n\n# 2) Rename the R variables inside the annotated source-code block.\n#    Each occurrence lives entirely inside its own run, so a plain\n#    find/replace keeps every other run (and its formatting) intact.\nfunction Replace-All($findText, $replaceText) {\n    $f = $d.Content.Find\n    $f.ClearFormatting()\n    $f.Replacement.ClearFormatting()\n    $f.Text = $findText\n    $f.Replacement.Text = $replaceText\n    $f.MatchCase = $true\n    $f.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2) | Out-Null\n}\n\nReplace-All \"p_0\" \"p0\"\nReplace-All \"se_0\" \"se0\"\n"}
